# ---------------------------------------------------------------------------
# Applies the "full country and no country dataframes correct" edit:
#   - params sheet: rename two variables (power_latop -> power_A,
#     time_laptop -> time_A), turn the old "energy_intensity_network" /
#     "bitrate_laptop" rows into new "power_B" / "time_B" interp rows,
#     clear the old "carbon_intensity" row (row 6) down to a single
#     formatted-but-empty date cell, and move the active selection to A6.
#   - rename the "energy_intensity_network" sheet to "time_B" and update
#     its per-country data + selection.
#   - duplicate that sheet into a brand-new "power_B" sheet (inserted right
#     after time_B, before "changes") with its own per-country data.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. params sheet -------------------------------------------------------
$params = $wb.Worksheets.Item("params")
$params.Activate() | Out-Null

# Row 2: power_latop -> power_A (rest of the row already matches target)
$params.Range("A2").Value = "power_A"

# Row 3: time_laptop -> time_A (rest of the row already matches target)
$params.Range("A3").Value = "time_A"

# Row 4: energy_intensity_network/exp -> power_B/interp
$params.Range("A4").Value = "power_B"
$params.Range("C4").Value = "interp"
$params.Range("D4").Value = "linear"
$params.Range("E4").Value = '{"2020-01-01":10, "2031-06-01":9.5}'
$params.Range("F4").Value = 0
$params.Range("G4").Value = 4
$params.Range("H4").Value = 0.05
$params.Range("I4").Value = 43617
$params.Range("J4").Value = "W"
$params.Range("P4").Value = "what does it mean? How do collect this info?"
$params.Range("Q4").Value = "x"
$params.Range("R4").Value = "power draw of laptop"
# S4 (id) stays 2, unchanged.

# Row 5: bitrate_laptop/exp -> time_B/interp
$params.Range("A5").Value = "time_B"
$params.Range("C5").Value = "interp"
$params.Range("D5").Value = "linear"
$params.Range("E5").Value = '{"2020-01-01":100, "2031-06-01":95}'
$params.Range("E5").ClearFormats()
# F5 stays 0, unchanged.
$params.Range("G5").Value = 5
# H5 stays 0.05, I5 stays 43617, unchanged.
$params.Range("J5").Value = "minute"
$params.Range("J5").ClearFormats()
$params.Range("Q5").Value = "x"
# S5 (id) stays 3, unchanged.

# Row 6: drop the old carbon_intensity row entirely, keep only the
# formatted (but empty) date cell in column I.
$params.Range("A6:H6").ClearContents()
$params.Range("I6").ClearContents()
$params.Range("J6:S6").ClearContents()

$params.Range("A6").Select() | Out-Null

# --- 2. energy_intensity_network -> time_B ---------------------------------
$timeB = $wb.Worksheets.Item("energy_intensity_network")
$timeB.Name = "time_B"
$timeB.Activate() | Out-Null

# Row 2 (UK)
$timeB.Range("C2").Value = '{"2020-01-01":100, "2031-06-01":95}'
$timeB.Range("D2").Value = 0
$timeB.Range("E2").Value = 5
$timeB.Range("F2").Value = 0.05
$timeB.Range("G2").Value = 4
$timeB.Range("G2").NumberFormat = "0"

# Row 3 (DE)
$timeB.Range("C3").Value = '{"2020-01-01":100, "2031-06-01":95}'
$timeB.Range("D3").Value = 1
$timeB.Range("E3").Value = 5
$timeB.Range("F3").Value = 0.05
$timeB.Range("G3").Value = 5

$timeB.Range("G2:G3").Select() | Out-Null

# --- 3. duplicate time_B -> power_B (new sheet) -----------------------------
$timeB.Copy([System.Reflection.Missing]::Value, $timeB)
$powerB = $wb.Worksheets.Item($timeB.Index + 1)
$powerB.Name = "power_B"
$powerB.Activate() | Out-Null

# Row 2 (UK)
$powerB.Range("C2").Value = '{"2020-01-01":10, "2031-06-01":9.5}'
$powerB.Range("D2").Value = 0
$powerB.Range("E2").Value = 4
$powerB.Range("F2").Value = 0.05
$powerB.Range("G2").Value = 2
$powerB.Range("G2").ClearFormats()

# Row 3 (DE)
$powerB.Range("C3").Value = '{"2020-01-01":10, "2031-06-01":9.5}'
$powerB.Range("D3").Value = 0
$powerB.Range("E3").Value = 4
$powerB.Range("F3").Value = 0.05
$powerB.Range("G3").Value = 3

$powerB.Range("C29").Select() | Out-Null

# --- 4. restore time_B as the active/selected sheet -------------------------
$timeB.Activate() | Out-Null

Write-Host "Edit complete"
